# Update Rizka - Manage Modules - Search
# Inserts a new "Admin - Manage Modules" worksheet right before the
# "Admin - Logout" sheet and populates it with the module search/condition
# table.

$wb = $excel.ActiveWorkbook

# Insert the new sheet immediately before "Admin - Logout" so the tab order
# becomes: Login, Quiz-Edit, Quiz-Delete, Manage Modules, Logout.
$logoutSheet = $wb.Worksheets.Item("Admin - Logout")
$ws = $wb.Worksheets.Add($logoutSheet)
$ws.Name = "Admin - Manage Modules"

# Header row
$ws.Range("A1").Value = "namaModul"
$ws.Range("B1").Value = "condition"

# Data rows
$ws.Range("A2").Value = "security testing !@#$%^&*()__+_)"
$ws.Range("B2").Value = "passed"

$ws.Range("A3").Value = "NEOP Teller Cash & PDC, dan FAB"
$ws.Range("B3").Value = "passed"

$ws.Range("A4").Value = "Karyawan baru pada fungsi AR"
$ws.Range("B4").Value = "failed"

# Trailing blank row (still present in the used range)
$ws.Range("A5").Value = ""
$ws.Range("B5").Value = ""

# Formatting: copy the bold / bordered header look already used elsewhere in
# the workbook (e.g. "Admin - Logout"!A1) onto the header row, then add the
# vertical centering that's unique to this sheet's header style.
$logoutSheet.Range("A1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$ws.Range("A1:B1").VerticalAlignment = -4108

# Data rows reuse the existing filled/bordered row styles from the other
# sheets in the workbook.
$loginSheet = $wb.Worksheets.Item("Admin - Login")
$loginSheet.Range("A2").Copy()
$ws.Range("A2:B3").PasteSpecial(-4122)

$loginSheet.Range("A4").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)

$ws.Columns.Item(1).ColumnWidth = 31.42578125

$ws.Range("H12").Select()
